# Annexe 3 - Programme hebdo de travail S_17
# Adds the "Mise en place / Prise en main de Trello" entries to row 21
# and a new data row (30th) for "Eval transactions" on row 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: fill in the previously-empty "Accompagnement" / "Ressources" cells
$ws.Range("E21").Value = "Mise en place"
$ws.Range("F21").Value = "Prise en main de Trello"

# Row 22: brand new data row (30 avril - Eval transactions)
$ws.Range("C22").Value = 30
$ws.Range("D22").Value = "Développer des composants d'accès aux données"
$ws.Range("E22").Value = "Eval transactions"
$ws.Range("F22").Value = "EVAL transactions"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0.5
$ws.Range("I22").Value = Get-Date -Year 2021 -Month 4 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("J22").Value = "eval_transactions.sql"
$ws.Range("K22").Value = "DWWM/S_17/eval_developper_des_composants_d acces_aux_donnees"
$ws.Range("L22").Value = "Github"

# Match the author's final active-cell selection
$ws.Range("G22").Select()
